$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parms")

# Update N row (row 17, columns B:H) to 500000
$ws.Range("B17:H17").Value = 500000

# Update the active selection on the sheet to A8
$ws.Activate()
$ws.Range("A8").Select()
